$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.034.97'
$ws.Range("E2").Value = '  -4.22%  '

$ws.Range("D3").Value = '1.968.04'
$ws.Range("E3").Value = '  -4.11%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.34'
$ws.Range("E5").Value = '  -3.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.627'
$ws.Range("E6").Value = '  -3.50%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.05'
$ws.Range("E7").Value = '  -3.76%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.370'
$ws.Range("E9").Value = '  -2.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '56.21'
$ws.Range("E10").Value = '  -5.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0816'
$ws.Range("E11").Value = '  +7.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.103'
$ws.Range("E12").Value = '  -1.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.862'
$ws.Range("E13").Value = '  -5.82%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.39'
$ws.Range("E14").Value = '  +7.31%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.14'
$ws.Range("E15").Value = '  -7.04%  '

$ws.Range("D16").Value = '2.254.71'
$ws.Range("E16").Value = '  -4.12%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.45'
$ws.Range("E17").Value = '  -3.28%  '

$ws.Range("D18").Value = '1.976.57'
$ws.Range("E18").Value = '  -3.67%  '

$ws.Range("D19").Value = '35.905.58'
$ws.Range("E19").Value = '  -4.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.10'
$ws.Range("E20").Value = '  -2.98%  '

$ws.Range("D21").Value = '0.0₃0858'
$ws.Range("E21").Value = '  -2.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '237.86'
$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.23'
$ws.Range("E23").Value = '  -2.66%  '

$ws.Range("E24").Value = '  -0.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.54'
$ws.Range("E25").Value = '  -8.52%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("E26").Value = '  -2.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.88'
$ws.Range("E27").Value = '  +2.56%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.00'
$ws.Range("E28").Value = '  -4.18%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.79'
$ws.Range("E29").Value = '  -0.85%  '

$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.131'
$ws.Range("E30").Value = '  +15.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.120'
$ws.Range("E31").Value = '  -1.57%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.90'
$ws.Range("E32").Value = '  -6.51%  '

$ws.Range("E33").Value = '  -7.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0623'
$ws.Range("E34").Value = '  +0.88%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.43'
$ws.Range("E35").Value = '  -6.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.32'
$ws.Range("E36").Value = '  +5.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.30'
$ws.Range("E37").Value = '  -6.82%  '

$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("E39").Value = '  +1.49%  '

$ws.Range("E40").Value = '  +14.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0989'
$ws.Range("E41").Value = '  -5.14%  '

$ws.Range("E42").Value = '  -0.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0212'
$ws.Range("E43").Value = '  -3.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.80'

$ws.Range("E45").Value = '  -5.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.29'
$ws.Range("E46").Value = '  -5.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '92.63'
$ws.Range("E47").Value = '  -3.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.58'
$ws.Range("E48").Value = '  -7.43%  '

$ws.Range("D49").Value = '1.344.47'
$ws.Range("E49").Value = '  -5.75%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.78'
$ws.Range("E50").Value = '  -5.73%  '

$ws.Range("D51").Value = '2.146.75'
$ws.Range("E51").Value = '  -3.95%  '
